$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "306.37"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.44%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.50"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-5.69%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.079"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.44%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07737"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-6.17%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.337"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.32%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.887"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-7.70%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "8.180"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.06%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-8.39%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9204"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.14%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1220"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-9.75%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1868"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-6.39%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08812"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-2.70%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03399"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.93%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09704"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.03%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001369"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.73%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006041"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.69%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.583"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.81%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-2.43%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1268"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-4.08%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.015"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.06%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.51%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "5,153.31%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04325"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.70%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001210"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.20%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004230"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-11.81%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.69%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02178"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-5.43%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04894"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-5.63%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007589"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.19%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009907"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-5.51%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1340"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.80%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002060"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.13%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009823"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.40%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006533"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-5.63%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.09%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "3.68%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.09%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.09%"
